$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the "Date: " / "__________" / "    Module: __________"
# runs into a single run with the combined text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Date: __________    Module: __________", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Date: __________    Module: __________", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: expand "The ratio of graphic elements and text is balanced
# and they are clearly separated." into the longer sentence, split
# across several runs (all sharing the original run's formatting).
# ---------------------------------------------------------------------
$rng = $d.Content
$old = "The ratio of graphic elements and text is balanced and they are clearly separated."
$rng.Find.Execute($old) | Out-Null

$xml = @"
<?xml version="1.0" encoding="utf-8"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="5502188A" w14:textId="68C9965A" w:rsidR="00C800E8" w:rsidRPr="00DF2B56" w:rsidRDefault="00DF2B56" w:rsidP="00C800E8">
            <w:pPr>
              <w:pStyle w:val="p1"/>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r w:rsidRPr="00DF2B56">
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">The ratio of graphic </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">and non-textual </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">elements </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">(formulas, tables, </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">lists, </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">etc.) </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">and </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">continuous </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue Light" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica Neue Light" w:cstheme="minorBidi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:val="en-US" w:eastAsia="en-US"/>
              </w:rPr>
              <w:t>text is balanced and they are clearly separated.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$rng.InsertXML($xml) | Out-Null
